{"js": "// Replace the ten-fold set of three-digit-by-one-digit multiplication\n// prompts throughout the document body with the newly generated values.\n// Each (old, new) pair below is applied with a case-sensitive, whole-\n// document search-and-replace; the old strings are unique in the\n// document, so this maps 1:1 onto the corresponding table cell.\nconst replacements = [\n  [\"491\u00d74=\", \"401\u00d72=\"],\n  [\"192\u00d76=\", \"948\u00d77=\"],\n  [\"770\u00d78=\", \"860\u00d72=\"],\n  [\"184\u00d76=\", \"683\u00d79=\"],\n  [\"945\u00d78=\", \"306\u00d75=\"],\n  [\"296\u00d72=\", \"805\u00d78=\"],\n  [\"906\u00d74=\", \"419\u00d74=\"],\n  [\"497\u00d74=\", \"365\u00d77=\"],\n  [\"243\u00d74=\", \"286\u00d73=\"],\n  [\"742\u00d76=\", \"900\u00d79=\"],\n  [\"111\u00d72=\", \"964\u00d73=\"],\n  [\"439\u00d76=\", \"381\u00d78=\"],\n  [\"961\u00d78=\", \"449\u00d78=\"],\n  [\"461\u00d72=\", \"256\u00d77=\"],\n  [\"329\u00d79=\", \"305\u00d78=\"],\n  [\"691\u00d76=\", \"873\u00d72=\"],\n  [\"204\u00d73=\", \"195\u00d76=\"],\n  [\"848\u00d77=\", \"893\u00d75=\"],\n  [\"332\u00d74=\", \"713\u00d77=\"],\n  [\"333\u00d78=\", \"821\u00d76=\"],\n  [\"216\u00d72=\", \"176\u00d78=\"],\n  [\"419\u00d75=\", \"899\u00d75=\"],\n  [\"417\u00d77=\", \"557\u00d79=\"],\n  [\"991\u00d74=\", \"477\u00d78=\"],\n  [\"905\u00d79=\", \"828\u00d74=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the three-digit-by-one-digit multiplication prompts throughout\n# the document body with the newly generated values. Each (old, new) pair\n# is unique within the document, so a simple Find/Replace per pair maps\n# 1:1 onto the corresponding table cell.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"491\u00d74=\", \"401\u00d72=\"),\n    @(\"192\u00d76=\", \"948\u00d77=\"),\n    @(\"770\u00d78=\", \"860\u00d72=\"),\n    @(\"184\u00d76=\", \"683\u00d79=\"),\n    @(\"945\u00d78=\", \"306\u00d75=\"),\n    @(\"296\u00d72=\", \"805\u00d78=\"),\n    @(\"906\u00d74=\", \"419\u00d74=\"),\n    @(\"497\u00d74=\", \"365\u00d77=\"),\n    @(\"243\u00d74=\", \"286\u00d73=\"),\n    @(\"742\u00d76=\", \"900\u00d79=\"),\n    @(\"111\u00d72=\", \"964\u00d73=\"),\n    @(\"439\u00d76=\", \"381\u00d78=\"),\n    @(\"961\u00d78=\", \"449\u00d78=\"),\n    @(\"461\u00d72=\", \"256\u00d77=\"),\n    @(\"329\u00d79=\", \"305\u00d78=\"),\n    @(\"691\u00d76=\", \"873\u00d72=\"),\n    @(\"204\u00d73=\", \"195\u00d76=\"),\n    @(\"848\u00d77=\", \"893\u00d75=\"),\n    @(\"332\u00d74=\", \"713\u00d77=\"),\n    @(\"333\u00d78=\", \"821\u00d76=\"),\n    @(\"216\u00d72=\", \"176\u00d78=\"),\n    @(\"419\u00d75=\", \"899\u00d75=\"),\n    @(\"417\u00d77=\", \"557\u00d79=\"),\n    @(\"991\u00d74=\", \"477\u00d78=\"),\n    @(\"905\u00d79=\", \"828\u00d74=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1          # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)  # wdReplaceAll\n}\n"}
